# Generate Report for Handoff
# Rows (Source File Name) whose handoff report was (re)generated this run.
$targetRows = @(7, 8, 10, 11, 12, 14)

$wb = $excel.ActiveWorkbook

# zh-cn: mark Priority as "ht" and bump the Latest Handoff Datetime.
$wsZh = $wb.Worksheets.Item("zh-cn")
foreach ($r in $targetRows) {
    $wsZh.Cells.Item($r, 5).Value = "ht"
    $wsZh.Cells.Item($r, 8).Value = "2016-08-13 08:22:55"
}

# de-de: mark Priority as "ht" and bump the Latest Handoff Datetime.
$wsDe = $wb.Worksheets.Item("de-de")
foreach ($r in $targetRows) {
    $wsDe.Cells.Item($r, 5).Value = "ht"
    $wsDe.Cells.Item($r, 8).Value = "2016-08-13 08:23:05"
}

# Overview: the "Latest HO Xliff Generate Date" column mirrors de-de's handoff datetime.
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $targetRows) {
    $wsOverview.Cells.Item($r, 7).Value = "2016-08-13 08:23:05"
}
